$d = $word.ActiveDocument

# Find the paragraph ending the main body text ("...codigo que os manipula.")
# and insert a page-break paragraph right after it (before the two trailing
# empty paragraphs).
$targetPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*alterar o c*digo que os manipula.*") {
        $targetPara = $p
    }
}

$pageBreakXml = @'
<?xml version="1.0" encoding="UTF-8"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
      <w:r>
        <w:br w:type="page"/>
      </w:r>
    </w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@

$r = $targetPara.Range
$r.Collapse(0)
$r.InsertXML($pageBreakXml)

# Find the first now-empty paragraph that follows (originally right after
# the body text) and insert the new "Sem projeto" section right after it,
# before the second empty paragraph / "Referencias" heading.
$emptyPara = $targetPara.Next().Next()

$contentXml = @'
<?xml version="1.0" encoding="UTF-8"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
      <w:r>
        <w:t xml:space="preserve">Sem projeto </w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>Como a maioria dos grandes aplicativos de grande porte tem o seu desenvolvimento de modo procedural, ocorre dificuldades na integração de novos módulos criados por diferentes setores e posteriormente é causado um grande custo para manutenções e atualizações (AGOSTINI; DECKER; SILVA, 2002).</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>Para um sistema de grande porte, é recomendado utilizar o paradigma de programação orientado a objetos para sua estrutura, por possibilitar modular, reutilizar e facilitar manutenções futuras (RAUT, 2020).</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t xml:space="preserve">Na Modelagem Orientada a Objetos, é priorizado a clareza e organização, usando abstrações do mundo real para entendimento, facilitando ao máximo a construção, manutenções e atualizações futuras. Este conceito de </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:bCs/>
        </w:rPr>
        <w:t>Abstração</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> possibilita a utilização da criação de um sistema feito com conceitos de objetos, que devem ser devidamente encapsulados (AGOSTINI; DECKER; SILVA, 2002).</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>Na abstração, os objetos são representações de entidades, conceitos ou processos do mundo real, focando em suas características essenciais e desconsiderando detalhes irrelevantes (AGOSTINI; DECKER; SILVA, 2002).</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t xml:space="preserve">Este conceito está diretamente ligado ao </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:bCs/>
        </w:rPr>
        <w:t>Encapsulamento</w:t>
      </w:r>
      <w:r>
        <w:t>, que consiste em agrupar os dados (atributos) e os métodos que os manipulam dentro de uma única unidade, a classe. A classe encapsula sua lógica interna, protegendo seus dados de acessos indevidos e garantindo a integridade do objeto (AGOSTINI; DECKER; SILVA, 2002). Por exemplo, a alteração de um atributo não é feita diretamente, mas através de um método público que pode conter regras de negócio e validações.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t xml:space="preserve">Para organizar e promover o reuso de código, utiliza-se o conceito de </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:bCs/>
        </w:rPr>
        <w:t>Herança</w:t>
      </w:r>
      <w:r>
        <w:t>. Este pilar permite a criação de uma classe base, que contém atributos e métodos comuns. A partir dela, classes mais específicas (classes-filhas) são derivadas, herdando todas as características da classe base e adicionando as suas próprias.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t xml:space="preserve">A herança abre caminho para o </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:bCs/>
        </w:rPr>
        <w:t>Polimorfismo</w:t>
      </w:r>
      <w:r>
        <w:t>, que é a capacidade de um objeto de uma classe derivada ser tratado como um objeto de sua classe base (RAUT, 2020). Isso permite que um método definido na classe base seja sobrescrito (</w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:i/>
          <w:iCs/>
        </w:rPr>
        <w:t>override</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t>) pelas classes filhas com suas próprias implementações.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>Dessa forma, o sistema pode invocar o método em um objeto do tipo da classe base, e o comportamento correto será executado em tempo de execução, dependendo de qual classe-filha a instância realmente pertence. Isso simplifica o código, tornando-o mais flexível e extensível, pois novos tipos de classes derivadas podem ser adicionados no futuro sem a necessidade de alterar o código que os manipula.</w:t>
      </w:r>
    </w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@

$r2 = $emptyPara.Range
$r2.Collapse(0)
$r2.InsertXML($contentXml)

Write-Output "DONE"
